$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Networking"
$ws.Range("B19").Value = "All the subtopics"
$ws.Range("C19").Value = "https://console.bluemix.net/catalog/?category=network"
$ws.Range("D19").Value = "networking"

$ws.Range("A21").Value = "storage"
$ws.Range("B21").Value = "All the subtopics"
$ws.Range("C21").Value = "https://console.bluemix.net/catalog/?category=storage"
$ws.Range("D21").Value = "storage"

$ws.Range("D22").Select()
